$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 16:46"

$ws.Range("B6").Value = 55528
$ws.Range("C6").Value = 672
$ws.Range("E6").Value = 54358
$ws.Range("G6").Value = 11
$ws.Range("H6").Value = 791

$ws.Range("E19").Value = 2579
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 28

$ws.Range("A28").Value = "Luxemburgo"
$ws.Range("B28").Value = 1333
$ws.Range("C28").Value = 234
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 1319
$ws.Range("F28").Value = 3
$ws.Range("H28").Value = 8

$ws.Range("A29").Value = "Irlanda"
$ws.Range("B29").Value = 1329
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = 1317
$ws.Range("F29").Value = 29
$ws.Range("H29").Value = 7

$ws.Range("A30").Value = "Japon"
$ws.Range("B30").Value = 1193
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 285
$ws.Range("E30").Value = 865
$ws.Range("F30").Value = 54
$ws.Range("G30").Value = 0
$ws.Range("H30").Value = 43

$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 1173
$ws.Range("C31").Value = 91
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = 1142
$ws.Range("F31").Value = 2
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 28

$ws.Range("A32").Value = "Chile"
$ws.Range("B32").Value = 1142
$ws.Range("C32").Value = 220
$ws.Range("D32").Value = 22
$ws.Range("E32").Value = 1117
$ws.Range("F32").Value = 7
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 3

$ws.Range("E36").Value = 806
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 14

$ws.Range("E50").Value = 434
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 8

$ws.Range("E57").Value = 395
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 1

$ws.Range("A68").Value = "Bulgaria"
$ws.Range("B68").Value = 242
$ws.Range("C68").Value = 24
$ws.Range("D68").Value = 4
$ws.Range("E68").Value = 235
$ws.Range("F68").Value = 8
$ws.Range("H68").Value = 3

$ws.Range("A69").Value = "Taiwan"
$ws.Range("B69").Value = 235
$ws.Range("C69").Value = 19
$ws.Range("D69").Value = 29
$ws.Range("E69").Value = 204
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 2

$ws.Range("A70").Value = "Hungria"
$ws.Range("B70").Value = 226
$ws.Range("C70").Value = 39
$ws.Range("D70").Value = 21
$ws.Range("E70").Value = 195
$ws.Range("F70").Value = 6
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 10

$ws.Range("A71").Value = "Letonia"
$ws.Range("B71").Value = 221
$ws.Range("C71").Value = 24
$ws.Range("D71").Value = 1
$ws.Range("E71").Value = 220
$ws.Range("F71").Value = 0
$ws.Range("H71").Value = 0

$ws.Range("B80").Value = 173
$ws.Range("C80").Value = 59
$ws.Range("E80").Value = 167

$ws.Range("A115").Value = "Consejo Danes para los Refugiados"
$ws.Range("C115").Value = 3
$ws.Range("F115").Value = 0

$ws.Range("A116").Value = "Mauricio"
$ws.Range("C116").Value = 6
$ws.Range("F116").Value = 1

$ws.Range("A129").Value = "Gibraltar"
$ws.Range("B129").Value = 26
$ws.Range("C129").Value = 11
$ws.Range("D129").Value = 5
$ws.Range("E129").Value = 21

$ws.Range("A130").Value = "Polinesia Francesa"
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 0
$ws.Range("E130").Value = 25
$ws.Range("H130").Value = 0

$ws.Range("A131").Value = "Jamaica"
$ws.Range("B131").Value = 25
$ws.Range("C131").Value = 4
$ws.Range("D131").Value = 2
$ws.Range("E131").Value = 22
$ws.Range("H131").Value = 1

$ws.Range("A132").Value = "Isla de Man"
$ws.Range("D132").Value = 0
$ws.Range("E132").Value = 23

$ws.Range("A133").Value = "Monaco"
$ws.Range("C133").Value = 0

$ws.Range("A134").Value = "Togo"
$ws.Range("C134").Value = 3
$ws.Range("D134").Value = 1
$ws.Range("E134").Value = 22

$ws.Range("A135").Value = "Guayana Francesa"
$ws.Range("B135").Value = 23
$ws.Range("D135").Value = 6
$ws.Range("E135").Value = 17
$ws.Range("H135").Value = 0

$ws.Range("A136").Value = "Guatemala"
$ws.Range("B136").Value = 21
$ws.Range("C136").Value = 0
$ws.Range("E136").Value = 20
$ws.Range("H136").Value = 1

$ws.Range("A137").Value = "Madagascar"
$ws.Range("B137").Value = 19
$ws.Range("C137").Value = 2
$ws.Range("E137").Value = 19

$ws.Range("A138").Value = "Barbados"
$ws.Range("B138").Value = 18
$ws.Range("E138").Value = 18

$ws.Range("A139").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("D139").Value = 0
$ws.Range("E139").Value = 17

$ws.Range("A140").Value = "Aruba"
$ws.Range("B140").Value = 17
$ws.Range("D140").Value = 1
$ws.Range("E140").Value = 16

$ws.Range("A149").Value = "El Salvador"
$ws.Range("C149").Value = 4

$ws.Range("A150").Value = "Guinea Ecuatorial"
$ws.Range("C150").Value = 0

$ws.Range("A154").Value = "Dominica"

$ws.Range("A155").Value = "Haiti"

$ws.Range("A158").Value = "Bermudas"

$ws.Range("A159").Value = "Benin"

$ws.Range("A160").Value = "Curazao"

$ws.Range("A162").Value = "Gabon"

$ws.Range("A164").Value = "Mozambique"
$ws.Range("C164").Value = 2
$ws.Range("D164").Value = 0
$ws.Range("E164").Value = 5

$ws.Range("A165").Value = "Bahamas"
$ws.Range("D165").Value = 1
$ws.Range("H165").Value = 0

$ws.Range("A166").Value = "Guyana"
$ws.Range("D166").Value = 0
$ws.Range("E166").Value = 4
$ws.Range("H166").Value = 1

$ws.Range("A167").Value = "Groenlandia"
$ws.Range("B167").Value = 5
$ws.Range("D167").Value = 2
$ws.Range("E167").Value = 3

$ws.Range("A168").Value = "Congo"

$ws.Range("A169").Value = "Suazilandia"

$ws.Range("A170").Value = "Guinea"

$ws.Range("A171").Value = "Santa Sede"
$ws.Range("C171").Value = 0
$ws.Range("E171").Value = 4
$ws.Range("H171").Value = 0

$ws.Range("A172").Value = "Siria"
$ws.Range("B172").Value = 4
$ws.Range("C172").Value = 3
$ws.Range("E172").Value = 4

$ws.Range("A173").Value = "Cabo Verde"
$ws.Range("B173").Value = 4
$ws.Range("C173").Value = 1
$ws.Range("H173").Value = 1

$ws.Range("A174").Value = "Antigua y Barbuda"

$ws.Range("A175").Value = "Angola"

$ws.Range("A176").Value = "Liberia"

$ws.Range("A177").Value = "Birmania"

$ws.Range("A178").Value = "San Bartolome"

$ws.Range("A179").Value = "Republica del Chad"

$ws.Range("A180").Value = "Santa Lucia"

$ws.Range("A182").Value = "Republica de Africa Central"
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 3

$ws.Range("A184").Value = "Zimbabue"

$ws.Range("A185").Value = "Nepal"
$ws.Range("C185").Value = 1
$ws.Range("D185").Value = 1
$ws.Range("H185").Value = 0

$ws.Range("A186").Value = "Sudan"
$ws.Range("B186").Value = 3
$ws.Range("H186").Value = 1

$ws.Range("A187").Value = "San Martin (Parte Holandesa)"
$ws.Range("C187").Value = 0

$ws.Range("A188").Value = "Guinea-Bisau"
$ws.Range("C188").Value = 2

$ws.Range("A190").Value = "Mali"
$ws.Range("C190").Value = 2

$ws.Range("A191").Value = "Nicaragua"
$ws.Range("C191").Value = 0

$ws.Range("A192").Value = "Butan"
$ws.Range("B192").Value = 2
$ws.Range("E192").Value = 2

$ws.Range("A193").Value = "Somalia"

$ws.Range("A194").Value = "Montserrat"

$ws.Range("A195").Value = "Eritrea"

$ws.Range("A196").Value = "Papua Nueva Guinea"

$ws.Range("A197").Value = "San Vicente y las Granadinas"

$ws.Range("A198").Value = "Islas Turcas y Caicos"

$ws.Range("A199").Value = "Timor Oriental"

$ws.Range("A200").Value = "Granada"

$ws.Range("A201").Value = "Libia"

$ws.Range("A202").Value = "Belice"
